$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the color values between the "stay" row (9) and "New Comers" row (10)
$stayColor = $ws.Range("C9").Value2
$newComersColor = $ws.Range("C10").Value2

$ws.Range("C9").Value = $newComersColor
$ws.Range("C10").Value = $stayColor

# Update the active selection to C10
$ws.Range("C10").Select()
